$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F13").Value = 63
$ws.Range("F14").Value = 80
$ws.Range("F15").Value = 77
$ws.Range("F16").Value = 66
$ws.Range("F17").Value = 62
$ws.Range("F18").Value = 59
$ws.Range("F19").Value = 74
$ws.Range("F20").Value = 59
$ws.Range("F21").Value = 64
$ws.Range("F22").Value = 74
$ws.Range("F23").Value = 68
$ws.Range("F24").Value = 72
$ws.Range("F25").Value = 78
$ws.Range("F26").Value = 52
$ws.Range("F27").Value = 72
$ws.Range("F28").Value = 64
$ws.Range("F29").Value = 67
$ws.Range("F30").Value = 60
$ws.Range("F31").Value = 66
$ws.Range("F32").Value = 72
$ws.Range("F33").Value = 73
$ws.Range("F34").Value = 80
$ws.Range("F35").Value = 60
$ws.Range("F36").Value = 64
$ws.Range("F37").Value = 77
$ws.Range("F38").Value = 76
$ws.Range("F39").Value = 76
$ws.Range("F40").Value = 80
$ws.Range("F46").Value = 56
$ws.Range("F47").Value = 97
$ws.Range("F48").Value = 50
$ws.Range("F49").Value = 67
$ws.Range("F50").Value = 51
$ws.Range("F51").Value = 62
$ws.Range("F52").Value = 54
$ws.Range("F53").Value = 89
$ws.Range("F54").Value = 81
$ws.Range("F55").Value = 70
$ws.Range("F56").Value = 92
$ws.Range("F57").Value = 59
$ws.Range("F58").Value = 55
$ws.Range("F59").Value = 81
$ws.Range("F60").Value = 64
$ws.Range("F61").Value = 69
$ws.Range("F62").Value = 63
$ws.Range("F63").Value = 93
$ws.Range("F64").Value = 53
$ws.Range("F65").Value = 86
$ws.Range("F66").Value = 93
$ws.Range("F67").Value = 98
$ws.Range("F68").Value = 79
$ws.Range("F69").Value = 70
$ws.Range("F70").Value = 55
$ws.Range("F71").Value = 62
$ws.Range("F72").Value = 69
$ws.Range("F73").Value = 82
$ws.Range("F79").Value = 87
$ws.Range("F80").Value = 52
$ws.Range("F81").Value = 90
$ws.Range("F82").Value = 60
$ws.Range("F83").Value = 58
$ws.Range("F84").Value = 75
$ws.Range("F85").Value = 67
$ws.Range("F86").Value = 86
$ws.Range("F87").Value = 80
$ws.Range("F88").Value = 87
$ws.Range("F89").Value = 53
$ws.Range("F90").Value = 75
$ws.Range("F91").Value = 71
$ws.Range("F92").Value = 74
$ws.Range("F93").Value = 67
$ws.Range("F94").Value = 66
$ws.Range("F95").Value = 67
$ws.Range("F96").Value = 56
$ws.Range("F97").Value = 65
$ws.Range("F98").Value = 56
$ws.Range("F99").Value = 82
$ws.Range("F100").Value = 52
$ws.Range("F101").Value = 95
$ws.Range("F102").Value = 73
$ws.Range("F103").Value = 99
$ws.Range("F104").Value = 80
$ws.Range("F105").Value = 99
$ws.Range("F106").Value = 78
$ws.Range("F112").Value = 61
$ws.Range("F113").Value = 98
$ws.Range("F114").Value = 71
$ws.Range("F115").Value = 86
$ws.Range("F116").Value = 84
$ws.Range("F117").Value = 98
$ws.Range("F118").Value = 90
$ws.Range("F119").Value = 96
$ws.Range("F120").Value = 80
$ws.Range("F121").Value = 95
$ws.Range("F122").Value = 95
$ws.Range("F123").Value = 95
$ws.Range("F124").Value = 94
$ws.Range("F125").Value = 75
$ws.Range("F126").Value = 67
$ws.Range("F127").Value = 80
$ws.Range("F128").Value = 72
$ws.Range("F129").Value = 100
$ws.Range("F130").Value = 54
$ws.Range("F131").Value = 86
$ws.Range("F132").Value = 93
$ws.Range("F133").Value = 84
$ws.Range("F134").Value = 55
$ws.Range("F135").Value = 91
$ws.Range("F136").Value = 50
$ws.Range("F137").Value = 53
$ws.Range("F138").Value = 52
$ws.Range("F139").Value = 55
$ws.Range("F145").Value = 75
$ws.Range("F146").Value = 63
$ws.Range("F147").Value = 87
$ws.Range("F148").Value = 55
$ws.Range("F149").Value = 63
$ws.Range("F150").Value = 53
$ws.Range("F151").Value = 84
$ws.Range("F152").Value = 53
$ws.Range("F153").Value = 88
$ws.Range("F154").Value = 95
$ws.Range("F155").Value = 51
$ws.Range("F156").Value = 61
$ws.Range("F157").Value = 76
$ws.Range("F158").Value = 68
$ws.Range("F159").Value = 98
$ws.Range("F160").Value = 81
$ws.Range("F161").Value = 55
$ws.Range("F162").Value = 55
$ws.Range("F163").Value = 96
$ws.Range("F164").Value = 99
$ws.Range("F165").Value = 99
$ws.Range("F166").Value = 99
$ws.Range("F167").Value = 95
$ws.Range("F168").Value = 81
$ws.Range("F169").Value = 78
$ws.Range("F170").Value = 74
$ws.Range("F171").Value = 73
$ws.Range("F172").Value = 53
$ws.Range("F178").Value = 77
$ws.Range("F179").Value = 89
$ws.Range("F180").Value = 90
$ws.Range("F181").Value = 97
$ws.Range("F182").Value = 87
$ws.Range("F183").Value = 52
$ws.Range("F184").Value = 70
$ws.Range("F185").Value = 78
$ws.Range("F186").Value = 87
$ws.Range("F187").Value = 84
$ws.Range("F188").Value = 97
$ws.Range("F189").Value = 51
$ws.Range("F190").Value = 90
$ws.Range("F191").Value = 84
$ws.Range("F192").Value = 71
$ws.Range("F193").Value = 56
$ws.Range("F194").Value = 89
$ws.Range("F195").Value = 76
$ws.Range("F196").Value = 70
$ws.Range("F197").Value = 71
$ws.Range("F198").Value = 69
$ws.Range("F199").Value = 65
$ws.Range("F200").Value = 78
$ws.Range("F201").Value = 77
$ws.Range("F202").Value = 57
$ws.Range("F203").Value = 80
$ws.Range("F204").Value = 58
$ws.Range("F205").Value = 92
